$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stats")

# Header for new column I, same style as the other header cells (A1:H1)
$ws.Range("I1").Value = "eta²"
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Effect size (eta squared) values for rows 2-10
$values = @(0.04, 0.06, 0.06, 0.16, 0.08, 0.03, 0.12, 0.2, 0.08)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i]
}
